# Apply cryptos list update (Sat Jan 27 19:53:47 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.949.14"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.275.88"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.68"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.22"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.82"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.70"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "2.628.95"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.37"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "2.278.07"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "41.867.93"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  +5.64%  "
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.04"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.91"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.03"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.32"
$ws.Range("E30").Value = "  +3.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.16"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.40"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0745"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.32"
$ws.Range("E36").Value = "  +4.13%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.66"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "2.008.59"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("E44").Value = "  +11.99%  "
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.26"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.63"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("E51").Value = "  +0.39%  "
